# For three of the four sensor-log worksheets, the tail of the log contains
# several "incomplete" timestamp groups (groups of rows that all share the
# same column-A timestamp, normally 6 rows per group, but truncated at the
# very end of the export). This edit appends one more duplicate row to each
# of the last 12 such groups (the last four incomplete group-triples),
# bringing their row counts from 7/5/3/1 up to 8/6/4/2, and shifting every
# row below each insertion point down by one.
#
# For every affected worksheet we insert a new row immediately below the
# last row of each of those 12 groups (processing from the bottom of the
# sheet upward so earlier insertion points keep their original row numbers)
# and fill the new row with a copy of the row directly above it.

$wb = $excel.ActiveWorkbook

# Map: worksheet index (1-based, matches workbook.xml sheet order) ->
# list of "insert a duplicate row right after this row" positions, given in
# the ORIGINAL (pre-edit) row numbering and ordered from bottom to top.
$targets = @{
    1 = @(517, 516, 515, 514, 511, 508, 505, 500, 495, 490, 483, 476)
    2 = @(531, 530, 529, 528, 525, 522, 519, 514, 509, 504, 497, 490)
    4 = @(535, 534, 533, 532, 529, 526, 523, 518, 513, 508, 501, 494)
}

foreach ($sheetIndex in $targets.Keys) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    $rowsToDuplicate = $targets[$sheetIndex]

    foreach ($srcRow in $rowsToDuplicate) {
        $destRow = $srcRow + 1

        # Push everything from destRow down by one row.
        $ws.Rows.Item($destRow).Insert()

        # Fill the freshly inserted (now-blank) row with a copy of the row
        # that used to sit right above the insertion point, values and
        # formatting both, so styling matches the surrounding rows exactly.
        $srcRange = $ws.Range("A" + $srcRow + ":I" + $srcRow)
        $destRange = $ws.Range("A" + $destRow + ":I" + $destRow)
        $srcRange.Copy($destRange)
    }
}
